# Fruta / hortaliza, semanal
# Insert a new weekly record for "Feria Lagunitas de Puerto Montt - Membrillo"
# at row 94, pushing all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(94).Insert()

$ws.Cells.Item(94, 1).Value  = 4
$ws.Cells.Item(94, 2).Value  = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(94, 3).Value  = 'Los Lagos'
$ws.Cells.Item(94, 4).Value  = [datetime]'2023-05-26'
$ws.Cells.Item(94, 5).Value  = 10
$ws.Cells.Item(94, 6).Value  = 'Fruta'
$ws.Cells.Item(94, 7).Value  = 100104
$ws.Cells.Item(94, 8).Value  = 'Frutos de pepita'
$ws.Cells.Item(94, 9).Value  = 100104003
$ws.Cells.Item(94, 10).Value = 'Membrillo'
$ws.Cells.Item(94, 11).Value = 'Champion'
$ws.Cells.Item(94, 12).Value = 'Primera'
$ws.Cells.Item(94, 13).Value = 200
$ws.Cells.Item(94, 14).Value = 13000
$ws.Cells.Item(94, 15).Value = 14000
$ws.Cells.Item(94, 16).Value = 13500
$ws.Cells.Item(94, 17).Value = '$/caja 18 kilos empedrada'
$ws.Cells.Item(94, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(94, 19).Value = 750
$ws.Cells.Item(94, 20).Value = 18
